$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 11:35"

# Row 15 - India
$ws.Range("B15").Value = 71339
$ws.Range("C15").Value = 571
$ws.Range("D15").Value = 23033
$ws.Range("E15").Value = 45996
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 2310

# Row 18 - Belgica
$ws.Range("B18").Value = 53779
$ws.Range("C18").Value = 330
$ws.Range("D18").Value = 13732
$ws.Range("E18").Value = 31286
$ws.Range("F18").Value = 465
$ws.Range("G18").Value = 54
$ws.Range("H18").Value = 8761

# Row 37 - Austria
$ws.Range("B37").Value = 15961
$ws.Range("C37").Value = 79
$ws.Range("D37").Value = 14148
$ws.Range("E37").Value = 1190
$ws.Range("F37").Value = 59
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 623

# Row 54 - Malasia
$ws.Range("B54").Value = 6742
$ws.Range("C54").Value = 16
$ws.Range("D54").Value = 5223
$ws.Range("E54").Value = 1410
$ws.Range("F54").Value = 16

# Row 91 - Eslovenia
$ws.Range("B91").Value = 1461
$ws.Range("C91").Value = 1
$ws.Range("D91").Value = 259
$ws.Range("E91").Value = 1100
$ws.Range("F91").Value = 9

# Row 95 - Hong Kong
$ws.Range("D95").Value = 991
$ws.Range("E95").Value = 53
